$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3 content updates -------------------------------------------------
$ws.Range("L3").Value = "Test agenda lagi"
$ws.Range("M3").Value = "notes"
$ws.Range("P3").Value = "2022-11-03"
$ws.Range("Q3").Formula = '=TEXT(TODAY(),"mmmm")'

# Touching WrapText (even to the value the cell already effectively has)
# makes Excel re-resolve the cell style record, dropping the redundant
# "applyFont" flag these cells carried (font is already the default font).
$ws.Range("L3").WrapText = $false
$ws.Range("N3").WrapText = $false
$ws.Range("O3:P3").WrapText = $false

# --- Sheet view: scroll position & active selection ------------------------
$ws.Application.ActiveWindow.ScrollColumn = 13
$ws.Range("X3").Select()
